# Final, renumbered member roster (post lend/return cleanup + admin delete/edit).
# Mirrors the app's "users" export: user_id, name, email, mobile, role.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$members = @(
    ,@("MEM-001","Admin","admin@library.com","","ADMIN")
    ,@("MEM-002","Sri","","","USER")
    ,@("MEM-003","Ganeshan","","","USER")
    ,@("MEM-004","Senthil","","","USER")
    ,@("MEM-005","muthu kumaresan","","","USER")
    ,@("MEM-006","Queretaro","","","USER")
    ,@("MEM-007","Kalai Priya","","","USER")
    ,@("MEM-008","Magesh (tamil class)","","","USER")
    ,@("MEM-009","Anand","","","USER")
    ,@("MEM-010","Gangesh","","","USER")
    ,@("MEM-011","Surya","","","USER")
    ,@("MEM-012","Selva","","","USER")
    ,@("MEM-013","Babu","","","USER")
    ,@("MEM-014","Shake","","","USER")
    ,@("MEM-015","Shanmugapriyan","","","USER")
    ,@("MEM-016","Hari","","","USER")
    ,@("MEM-017","Saidai Jagan","","","USER")
    ,@("MEM-018","benazir (tamil class)","","","USER")
    ,@("MEM-019","Muthusamy","","","USER")
    ,@("MEM-020","Sasikumar","","","USER")
    ,@("MEM-021","Maheswari","","","USER")
    ,@("MEM-022","Chandra prakash","","","USER")
    ,@("MEM-023","Naveen","","","USER")
    ,@("MEM-024","Muruganantham","","","USER")
    ,@("MEM-025","Vignesh (AI)","","","USER")
    ,@("MEM-026","Saidai Jagan","","","USER")
    ,@("MEM-027","Hariprasad","","","USER")
    ,@("MEM-028","Arun Bharani","","","USER")
    ,@("MEM-029","Suganya","","","USER")
    ,@("MEM-030","Priyanka sathish","","","USER")
    ,@("MEM-031","Susidevi","","","USER")
    ,@("MEM-032","Kasirajan","","","USER")
    ,@("MEM-033","varatharaj","","","USER")
    ,@("MEM-034","Hari Prakash","","","USER")
    ,@("MEM-035","Logeshwari","","","USER")
    ,@("MEM-036","Stela","","","USER")
    ,@("MEM-037","Mahesh","","","USER")
    ,@("MEM-038","SuryaDevi","","","USER")
    ,@("MEM-039","Vignesh","","","USER")
    ,@("MEM-040","sibi","","","USER")
    ,@("MEM-041","hariprasad","","","USER")
    ,@("MEM-042","Anbuvel","","","USER")
    ,@("MEM-043","Logeshwari","","","USER")
    ,@("MEM-044","Thoufik","","","USER")
    ,@("MEM-045","Divya","","","USER")
    ,@("MEM-046","Rajeswari","","","USER")
    ,@("MEM-047","Indrasenan","","","USER")
    ,@("MEM-048","sri","","","USER")
    ,@("MEM-049","banupriya","","","USER")
    ,@("MEM-050","Ragavan","","","USER")
    ,@("MEM-051","Vignesh V","vignesh29498@gmail.com",9080314434,"USER")
)

# Rewrite the member rows (row 2 onward) with the renumbered roster.
$row = 2
foreach ($m in $members) {
    $ws.Cells.Item($row, 1).Value = $m[0]          # user_id (re-sequenced MEM-xxx)
    $ws.Cells.Item($row, 2).Value = $m[1]          # name
    $ws.Cells.Item($row, 3).Value = $m[2]          # email (blank unless set)
    if ($m[3] -eq "") {
        $ws.Cells.Item($row, 4).Value = ""
    } else {
        $ws.Cells.Item($row, 4).Value = [double]$m[3]   # mobile (numeric)
    }
    $ws.Cells.Item($row, 5).Value = $m[4]          # role
    $row++
}

# Drop the now-unused trailing rows - roster shrank from 61 to 51 members.
$lastOldRow = 62
if ($lastOldRow -ge $row) {
    $ws.Range("A" + $row + ":E" + $lastOldRow).EntireRow.Delete()
}

Write-Output "Roster rewritten: $($members.Count) members (rows 2..$($row-1))."
